# Auto-sync batch FINAL (AUTO-TIMEOUT)
# Appends the newest batch of MFS/mobile-money ledger rows (86-107) to
# Sheet1, right after the existing data, and widens the sheet's used
# range accordingly.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$rows = @(
    @{ Row = 86;  Date = "2026-02-11 12:02:52"; Number = "237671823369"; Name = "MFS ENTREE COLLEGE MALANGUE";                                              Balance = 103700 },
    @{ Row = 87;  Date = "2026-02-11 13:50:43"; Number = "237672128028"; Name = "DELVIN NDIFON BAH";                                                         Balance = 64935 },
    @{ Row = 88;  Date = "2026-02-10 14:37:59"; Number = "237672277367"; Name = "TOP MOBIL KM5 LTDLA_POLAS_BTQ_KM5";                                         Balance = 1100121 },
    @{ Row = 89;  Date = "2026-02-11 15:23:12"; Number = "237674853971"; Name = "NJOSSEU TCHOUNZOU TOP MOBILE";                                              Balance = 379307 },
    @{ Row = 90;  Date = "2026-02-11 15:01:40"; Number = "237674884705"; Name = "BAH AMADOU MOUNTAGHA ETS MOBILE FINANCIAL SERVICES MFS";                    Balance = 19074 },
    @{ Row = 91;  Date = "2026-02-11 15:36:14"; Number = "237675779272"; Name = "RODES NGWEM KEMAYOU";                                                       Balance = 137655 },
    @{ Row = 92;  Date = "2026-02-11 14:54:38"; Number = "237677304210"; Name = "CARINE OROCK";                                                              Balance = 90031 },
    @{ Row = 93;  Date = "2026-02-06 09:50:11"; Number = "237678267353"; Name = "LA NEGRESSE SARL EMBOLA BELTUS MBU";                                        Balance = 0 },
    @{ Row = 94;  Date = "2026-02-11 14:24:33"; Number = "237678370615"; Name = "ESSEN ONGOLONG BERTHE HORTENSE ETS MOBILE FINANCIAL SERVICES MFS";          Balance = 123243 },
    @{ Row = 95;  Date = "2026-02-11 15:14:13"; Number = "237678836319"; Name = "KAMDOM DOMINIQUE STEPHANIE ETS MOBILE FINANCIAL SERVICES MFS";              Balance = 109522 },
    @{ Row = 96;  Date = "2026-02-11 14:58:14"; Number = "237678922502"; Name = "NWOAGA TCHAMDJOU EPSE KAMSEU EMILINE ETS LE CONTENT";                       Balance = 563253 },
    @{ Row = 97;  Date = "2026-02-11 14:10:23"; Number = "237679884264"; Name = "MFS CICAM";                                                                 Balance = 117365 },
    @{ Row = 98;  Date = "2026-02-11 13:32:34"; Number = "237681019523"; Name = "ETS MOULAY RIPERT AND COMPANY";                                             Balance = 90676 },
    @{ Row = 99;  Date = "2026-02-11 14:54:37"; Number = "237681125655"; Name = "EMILE MADELO";                                                              Balance = 12821 },
    @{ Row = 100; Date = "2026-02-11 15:09:32"; Number = "237681240793"; Name = "MBANE EMILIE FRANCOISE ETS MOBILE FINANCIAL SERVICES MFS";                  Balance = 2959 },
    @{ Row = 101; Date = "2026-02-11 12:38:08"; Number = "237682117915"; Name = "MEKUEKO FOUDJO BERLINE DIDIANE ETS MOBILE FINANCIAL SERVICES MFS";          Balance = 23687 },
    @{ Row = 102; Date = "2026-02-11 14:37:09"; Number = "237682154553"; Name = "N A ISUFUH MIEMONA NGESSY ETS MOBILE FINANCIAL SERVICES MFS";               Balance = 31551 },
    @{ Row = 103; Date = "2026-02-11 15:36:22"; Number = "237682803277"; Name = "NGUEMASSOM RENE MARTIAL LA NEGRESSE SARL";                                  Balance = 69625 },
    @{ Row = 104; Date = "2026-02-11 15:06:30"; Number = "237683323481"; Name = "ETS LE CONTENT TSAZE DONFOUET FLORETTE ROSINE";                             Balance = 509243 },
    @{ Row = 105; Date = "2026-02-11 15:29:21"; Number = "237683368985"; Name = "MFS BELL HENRIE BERNARD";                                                   Balance = 232132 },
    @{ Row = 106; Date = "2026-02-11 15:43:45"; Number = "237683432110"; Name = "ERODINE TOUMENI";                                                           Balance = 822255 },
    @{ Row = 107; Date = "2026-02-11 15:34:26"; Number = "237683743490"; Name = "ETS LE CONTENT NGAH MARIE";                                                 Balance = 4823 }
)

foreach ($r in $rows) {
    $rowIndex = $r.Row

    # Column A: free-text timestamp (kept as text, same as every row above it).
    $ws.Range("A$rowIndex").Value = $r.Date

    # Column B: the subscriber number is long digits-only text. Assigning a
    # plain numeric-looking string to a General-formatted cell would get
    # auto-coerced into a number, so mark the cell as Text first, write the
    # value, then drop the cell back onto the sheet's normal/default style
    # so no stray formatting is left behind on the new rows.
    $ws.Range("B$rowIndex").NumberFormat = "@"
    $ws.Range("B$rowIndex").Value = $r.Number
    $ws.Range("B$rowIndex").Style = "Normal"

    # Column C: merchant / account name, plain text.
    $ws.Range("C$rowIndex").Value = $r.Name

    # Column D: balance, numeric.
    $ws.Range("D$rowIndex").Value = $r.Balance
}
